$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(40, 73, "PT271224160107", "wwewe", 22, "", "", 45653.66744212963, 0, "Severe"),
  @(41, 74, "PT271224161405", "sss", 24, "", "", 45653.6764699074, 0, "Severe"),
  @(42, 75, "PT271224161438", "sss", 24, "", "", 45653.67685185185, 0, "Severe"),
  @(43, 76, "PT271224162553", "sdsad", 23, "", "", 45653.6846412037, 0, "Severe"),
  @(44, 77, "PT271224162714", "qqqwee", 12, "", "", 45653.68557870371, 0, "Severe"),
  @(45, 78, "PT271224171633", "www", 23, "", "", 45653.71982638889, 0, "Severe"),
  @(46, 79, "PT271224171848", "qeqqdsdd", 33, "", "", 45653.7216550926, 2, "Severe"),
  @(47, 82, "PT271224172242", "qqq", 33, "", "", 45653.72439814815, 5, "Severe"),
  @(48, 85, "PT271224174848", "wer", 23, "", "", 45653.74231481482, 0, "Severe"),
  @(49, 87, "PT271224174028", "ff", 12, "", "", 45653.74271990741, 3, "Severe"),
  @(50, 88, "PT271224175343", "ffff", 22, "", "", 45653.74563657407, 0, "Severe"),
  @(51, 89, "PT271224180528", "DGP", 25, "dementia", "Normal", 45653.75569444444, 28, "May be Normal"),
  @(52, 92, "PT271224181048", "sos", 88, "", "", 45653.7575, 0, "Severe")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    if ([string]::IsNullOrEmpty($row[5]) -eq $false) {
        $ws.Cells.Item($r, 5).Value = $row[5]
    }
    if ([string]::IsNullOrEmpty($row[6]) -eq $false) {
        $ws.Cells.Item($r, 6).Value = $row[6]
    }
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
}

